# Updated cryptos list (refreshed Price / Volume(1h) figures, plus a rank
# swap between THORChain and WEMIXToken at rows 36-37) as published by the
# "Updated cryptos list ... with GitHub Actions" scraper run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($rng, $val)
    # Force the cell to be treated as text even when the value looks numeric
    # (matches the source data which stores these as literal text strings).
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '36.400.20'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '1.932.42'
$ws.Range("E3").Value = '  -2.05%  '
$ws.Range("E4").Value = '  +0.01%  '
Set-TextValue $ws.Range("D5") '241.63'
$ws.Range("E5").Value = '  -1.34%  '
Set-TextValue $ws.Range("D6") '0.607'
$ws.Range("E6").Value = '  -2.46%  '
$ws.Range("E7").Value = '  +0.02%  '
Set-TextValue $ws.Range("D8") '56.43'
$ws.Range("E8").Value = '  -2.89%  '
$ws.Range("E9").Value = '  -2.89%  '
$ws.Range("E10").Value = '  +0.14%  '
$ws.Range("E11").Value = '  -1.98%  '
$ws.Range("D12").Value = '2.216.92'
$ws.Range("E12").Value = '  -1.95%  '
Set-TextValue $ws.Range("D13") '21.02'
$ws.Range("E13").Value = '  -6.79%  '
Set-TextValue $ws.Range("D14") '0.800'
$ws.Range("E14").Value = '  -6.07%  '
Set-TextValue $ws.Range("D15") '13.30'
$ws.Range("E15").Value = '  -3.82%  '
Set-TextValue $ws.Range("D16") '5.10'
$ws.Range("E16").Value = '  -5.76%  '
$ws.Range("D17").Value = '1.932.61'
$ws.Range("E17").Value = '  -2.15%  '
$ws.Range("D18").Value = '36.344.94'
$ws.Range("E18").Value = '  +0.49%  '
Set-TextValue $ws.Range("D19") '68.73'
$ws.Range("E19").Value = '  -2.28%  '
$ws.Range("D20").Value = '0.0₃0857'
$ws.Range("E20").Value = '  -2.54%  '
Set-TextValue $ws.Range("D21") '226.19'
$ws.Range("E21").Value = '  -2.99%  '
$ws.Range("E22").Value = '  -5.67%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("E24").Value = '  -7.45%  '
$ws.Range("E25").Value = '  -1.81%  '
Set-TextValue $ws.Range("D26") '9.02'
$ws.Range("E26").Value = '  -8.06%  '
Set-TextValue $ws.Range("D27") '160.75'
$ws.Range("E27").Value = '  -1.82%  '
Set-TextValue $ws.Range("D28") '0.134'
$ws.Range("E28").Value = '  +0.68%  '
Set-TextValue $ws.Range("D29") '19.12'
$ws.Range("E29").Value = '  -2.96%  '
$ws.Range("E30").Value = '  -1.63%  '
Set-TextValue $ws.Range("D31") '1.11'
$ws.Range("E31").Value = '  -5.67%  '
$ws.Range("E32").Value = '  -6.75%  '
Set-TextValue $ws.Range("D33") '0.0616'
$ws.Range("E33").Value = '  -6.42%  '
Set-TextValue $ws.Range("D34") '4.12'
$ws.Range("E34").Value = '  -7.14%  '
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D36") '1.78'
$ws.Range("E36").Value = '  -1.80%  '
$ws.Range("B37").Value = 'THORChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range("D37") '5.91'
$ws.Range("E37").Value = '  -3.25%  '
$ws.Range("E38").Value = '  -2.36%  '
Set-TextValue $ws.Range("D39") '3.02'
$ws.Range("E39").Value = '  +3.20%  '
Set-TextValue $ws.Range("D40") '0.0970'
$ws.Range("E40").Value = '  +1.17%  '
$ws.Range("E41").Value = '  -0.38%  '
$ws.Range("E42").Value = '  -2.21%  '
$ws.Range("E43").Value = '  -5.86%  '
Set-TextValue $ws.Range("D44") '15.44'
$ws.Range("E44").Value = '  -3.33%  '
$ws.Range("D45").Value = '1.332.72'
$ws.Range("E45").Value = '  -2.56%  '
$ws.Range("E46").Value = '  -6.53%  '
Set-TextValue $ws.Range("D47") '85.04'
$ws.Range("E47").Value = '  -6.43%  '
Set-TextValue $ws.Range("D48") '7.06'
$ws.Range("E48").Value = '  -4.21%  '
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("D50").Value = '2.109.50'
$ws.Range("E50").Value = '  -1.88%  '
Set-TextValue $ws.Range("D51") '43.22'
$ws.Range("E51").Value = '  -3.66%  '

Write-Output "applied 86 cell updates"
